# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# sheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 7005
$wsExhibit.Range("F5").Value = 458
$wsExhibit.Range("F7").Value = 6920
$wsExhibit.Range("F23").Value = 704

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 7005
$wsAll.Range("F5").Value = 458
$wsAll.Range("F7").Value = 6920
$wsAll.Range("F25").Value = 704
